# Consolidate the three text runs ("Below", " ", "section-level") in the
# title placeholder of slide 3 into a single run, matching the behaviour
# of a PowerPoint writer that merges adjacent runs when possible.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

$full = $tr.Text
$whole = $tr.Characters(1, $tr.Length)
$whole.Text = $full
